$d = $word.ActiveDocument

# The existing "_GoBack" bookmark sits right at the end of the paragraph
# that ends in "...srollbarral).". We split that paragraph there to add a
# new bold list item ("Jobb optimalizálás telefonos felülethez is (nincs
# teljesen kész.)") and recreate the bookmark in its new location, inside
# the freshly inserted paragraph, right before the closing ")".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("srollbarral).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text (right after "srollbarral).")
    # and split the paragraph there.
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(1, 1)

    $rng.InsertAfter("Jobb optimalizálás telefonos felülethez is")
    $rng.Collapse(0)

    $rng.InsertAfter(" (nincs teljesen kész.")
    $rng.Collapse(0)

    $rng.InsertAfter(")")

    # Re-create "_GoBack" collapsed right before the final ")" we just typed.
    $bmRange = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
